# B6-PowerPoint.pptx edit: Sat, Jul 11, 2020  1:05:13 PM
#
# 1) Re-colour the deck's theme (currently the "Integral" / "Red Violet"
#    colour scheme) to the standard Office colour scheme, as seen through
#    the Slide Master's theme (ppt/theme/theme2.xml - the theme actually
#    driving every slide in the deck).
# 2) Re-style the three data tables (slides 14, 15, 16) with a different
#    built-in table style.

$p = $ppt.ActivePresentation

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Theme colours: Red Violet -> Office -------------------------------
# ThemeColorScheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt($officeColors[$i - 1])
}

# --- 2) Table styles on slides 14, 15 and 16 -------------------------------
$newTableStyleId = "{83912559-4242-44FF-86F9-07726B5F9A10}"
foreach ($slideIndex in 14, 15, 16) {
    $tableShape = $p.Slides.Item($slideIndex).Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyleId, $true)
    }
}
